$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows that correspond to genomes / taxa removed by this commit
# (kept rows are sequences that have multiple PPO domains).
# Rows are deleted from bottom to top so row numbers stay valid.
$rowsToDelete = @(103, 95, 89, 84, 83, 82)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
